$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New row 71: "452. Minimum Number of Arrows to Burst Balloons"
# ---------------------------------------------------------------------------

$question = '452. Minimum Number of Arrows to Burst Balloons'
$difficulty = 'Medium'
$pattern = 'Intervals'
$notes = 'Greedy approach. First sort the balloons by their ending position with Integer.compare(a[1],b[1]), then iterate, shooting the arrow at the ending position of the current balloon, then start again from the next one. You iterate by points[i][1] and set the arrow there, then check if points[i][0] falls in that range. If it is, you continue without iterating arrow count or shifting the arrow position. For large values, you have to use Integer compare, or convert to long first.'
$link = 'https://leetcode.com/problems/minimum-number-of-arrows-to-burst-balloons/solutions/1686627/c-java-python-6-lines-sort-and-greedy-image-explanation/?envType=study-plan-v2&envId=leetcode-75 '

$ws.Range("A71").Value = $question
$ws.Range("B71").Value = $difficulty
$ws.Range("C71").Value = $pattern
$ws.Range("D71").Value = $notes
$ws.Range("E71").Value = $link

# Match the "Medium" fill formatting already used on column B (e.g. B70)
$ws.Range("B71").Interior.Color = $ws.Range("B70").Interior.Color()

# Turn E71 into a real hyperlink pointing at the (trimmed) URL, then re-apply
# the built-in Hyperlink style so the cell matches the rest of column E.
$ws.Hyperlinks.Add($ws.Range("E71"), "https://leetcode.com/problems/minimum-number-of-arrows-to-burst-balloons/solutions/1686627/c-java-python-6-lines-sort-and-greedy-image-explanation/?envType=study-plan-v2&envId=leetcode-75") | Out-Null
$ws.Range("E71").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Update the saved view state (best effort)
# ---------------------------------------------------------------------------

$ws.Activate()
$ws.Range("B46").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 2
$ws.Range("E74").Select() | Out-Null
